$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C4"  = 2.208165160720954
    "E4"  = 1.903751357432193
    "C5"  = 1.614140618728332
    "E5"  = 1.770808585446004
    "C6"  = 1.625793900975747
    "E6"  = 1.586821460965226
    "C7"  = 0.8049364973309325
    "E7"  = 1.421244400332
    "C8"  = 0.5775251578155283
    "E8"  = 1.341244385861273
    "C9"  = 1.901826580533572
    "E9"  = 1.53605963063923
    "C10" = 2.590339257583607
    "E10" = 1.672072534917302
    "C11" = 1.713587272940131
    "E11" = 1.721854626734953
    "C12" = 1.05432456490544
    "E12" = 1.415552619392124
    "C13" = 1.566023898188384
    "E13" = 1.644188696416427
    "C14" = 2.155932165770968
    "E14" = 1.805141163113122
    "C15" = 2.443967114785739
    "E15" = 2.026008136667135
    "C16" = 0.388123216496683
    "E16" = 1.819907598678561
    "C17" = -2.811030211656218
    "E17" = 0.8407670860975047
    "C18" = 1.250641979737566
    "E18" = 1.466559393695466
    "C19" = 2.302179720973463
    "E19" = 1.805984941845473
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
